$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 117.625
$ws.Range("I33").Value = 68.5
$ws.Range("J33").Value = 265
$ws.Range("K33").Value = 68.5
$ws.Range("L33").Value = 265
$ws.Range("M33").Value = 160.5
$ws.Range("N33").Value = -723
$ws.Range("H86").Value = 1764.4286
$ws.Range("I86").Value = 1658.5
$ws.Range("J86").Value = 2400
$ws.Range("K86").Value = 1658.5
$ws.Range("L86").Value = 2400
$ws.Range("M86").Value = -535.5
$ws.Range("N86").Value = -4646
$ws.Range("H89").Value = 1764.4286
$ws.Range("I89").Value = 1658.5
$ws.Range("J89").Value = 2400
$ws.Range("K89").Value = 8292.5
$ws.Range("L89").Value = 12000
$ws.Range("M89").Value = -2676.5
$ws.Range("N89").Value = -23232
$ws.Range("H92").Value = 556452.4399999999
$ws.Range("I92").Value = 654562
$ws.Range("J92").Value = 498.33334
$ws.Range("K92").Value = 654562
$ws.Range("L92").Value = 498.33334
$ws.Range("M92").Value = -653314
$ws.Range("N92").Value = -2994.33334
$ws.Range("H98").Value = 623054.1
$ws.Range("I98").Value = 861280.3
$ws.Range("J98").Value = 3666
$ws.Range("K98").Value = 861280.3
$ws.Range("L98").Value = 3666
$ws.Range("M98").Value = -859782.3
$ws.Range("N98").Value = -6662
$ws.Range("H106").Value = 12459953
$ws.Range("I106").Value = 14016575
$ws.Range("J106").Value = 6980
$ws.Range("K106").Value = 14016575
$ws.Range("L106").Value = 6980
$ws.Range("M106").Value = -14015944
$ws.Range("N106").Value = -8242
$ws.Range("H122").Value = 623054.1
$ws.Range("I122").Value = 861280.3
$ws.Range("J122").Value = 3666
$ws.Range("K122").Value = 2583840.9
$ws.Range("L122").Value = 10998
$ws.Range("M122").Value = -2581390.9
$ws.Range("N122").Value = -15898
$ws.Range("H132").Value = 316349.12
$ws.Range("I132").Value = 393040.75
$ws.Range("J132").Value = 52189.11
$ws.Range("K132").Value = 1179122.25
$ws.Range("L132").Value = 156567.33
$ws.Range("M132").Value = -1176592.25
$ws.Range("N132").Value = -161627.33
$ws.Range("H137").Value = 35715932
$ws.Range("I137").Value = 50000996
$ws.Range("J137").Value = 3273.75
$ws.Range("K137").Value = 150002988
$ws.Range("L137").Value = 9821.25
$ws.Range("M137").Value = -150000438
$ws.Range("N137").Value = -14921.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3618.9583
$ws.Range("I61").Value = 2357.5386
$ws.Range("K61").Value = 2357.5386
$ws.Range("M61").Value = -2145.5386
$ws.Range("H74").Value = 6585.174
$ws.Range("I74").Value = 1485.0588
$ws.Range("J74").Value = 21035.5
$ws.Range("K74").Value = 1485.0588
$ws.Range("L74").Value = 21035.5
$ws.Range("M74").Value = -611.0588
$ws.Range("N74").Value = -22783.5
$ws.Range("H77").Value = 6585.174
$ws.Range("I77").Value = 1485.0588
$ws.Range("J77").Value = 21035.5
$ws.Range("K77").Value = 7425.294
$ws.Range("L77").Value = 105177.5
$ws.Range("M77").Value = -3057.294
$ws.Range("N77").Value = -113913.5
$ws.Range("H132").Value = 2398.509
$ws.Range("I132").Value = 1754.2941
$ws.Range("J132").Value = 3441.524
$ws.Range("K132").Value = 5262.8823
$ws.Range("L132").Value = 10324.572
$ws.Range("M132").Value = -2732.8823
$ws.Range("N132").Value = -15384.572
$ws.Range("H136").Value = 3618.9583
$ws.Range("I136").Value = 2357.5386
$ws.Range("K136").Value = 7072.6158
$ws.Range("M136").Value = -4522.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6496.2104
$ws.Range("I86").Value = 1281.9286
$ws.Range("J86").Value = 9537.875
$ws.Range("K86").Value = 1281.9286
$ws.Range("L86").Value = 9537.875
$ws.Range("M86").Value = -158.9286
$ws.Range("N86").Value = -11783.875
$ws.Range("H89").Value = 6496.2104
$ws.Range("I89").Value = 1281.9286
$ws.Range("J89").Value = 9537.875
$ws.Range("K89").Value = 6409.643
$ws.Range("L89").Value = 47689.375
$ws.Range("M89").Value = -793.643
$ws.Range("N89").Value = -58921.375
$ws.Range("H94").Value = 754.6818
$ws.Range("I94").Value = 689.0625
$ws.Range("K94").Value = 689.0625
$ws.Range("M94").Value = -238.0625
$ws.Range("H134").Value = 29415260
$ws.Range("I134").Value = 50002050
$ws.Range("K134").Value = 150006150
$ws.Range("M134").Value = -150003615

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2092.7856
$ws.Range("I31").Value = 1950
$ws.Range("J31").Value = 2449.75
$ws.Range("K31").Value = 1950
$ws.Range("L31").Value = 2449.75
$ws.Range("M31").Value = -1655
$ws.Range("N31").Value = -3039.75
$ws.Range("H34").Value = 2092.7856
$ws.Range("I34").Value = 1950
$ws.Range("J34").Value = 2449.75
$ws.Range("K34").Value = 1950
$ws.Range("L34").Value = 2449.75
$ws.Range("M34").Value = -1748
$ws.Range("N34").Value = -2853.75
$ws.Range("H58").Value = 2006.3243
$ws.Range("I58").Value = 1270.9615
$ws.Range("J58").Value = 3744.4546
$ws.Range("K58").Value = 1270.9615
$ws.Range("L58").Value = 3744.4546
$ws.Range("M58").Value = -1067.9615
$ws.Range("N58").Value = -4150.4546
$ws.Range("H64").Value = 19500
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 19500
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H132").Value = 2757.1
$ws.Range("I132").Value = 2292.7368
$ws.Range("J132").Value = 3559.182
$ws.Range("K132").Value = 6878.2104
$ws.Range("L132").Value = 10677.546
$ws.Range("M132").Value = -4348.2104
$ws.Range("N132").Value = -15737.546
$ws.Range("H134").Value = 2919.9565
$ws.Range("I134").Value = 1368.1818
$ws.Range("J134").Value = 4342.4165
$ws.Range("K134").Value = 4104.5454
$ws.Range("L134").Value = 13027.2495
$ws.Range("M134").Value = -1569.5454
$ws.Range("N134").Value = -18097.2495
$ws.Range("H136").Value = 2006.3243
$ws.Range("I136").Value = 1270.9615
$ws.Range("J136").Value = 3744.4546
$ws.Range("K136").Value = 3812.8845
$ws.Range("L136").Value = 11233.3638
$ws.Range("M136").Value = -1262.8845
$ws.Range("N136").Value = -16333.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2731
$ws.Range("I80").Value = 2640
$ws.Range("J80").Value = 2776.5
$ws.Range("K80").Value = 2640
$ws.Range("L80").Value = 2776.5
$ws.Range("M80").Value = -1642
$ws.Range("N80").Value = -4772.5
$ws.Range("H83").Value = 2731
$ws.Range("I83").Value = 2640
$ws.Range("J83").Value = 2776.5
$ws.Range("K83").Value = 13200
$ws.Range("L83").Value = 13882.5
$ws.Range("M83").Value = -8208
$ws.Range("N83").Value = -23866.5
$ws.Range("H132").Value = 2793.3555
$ws.Range("I132").Value = 2583.423
$ws.Range("J132").Value = 3080.6316
$ws.Range("K132").Value = 7750.268999999999
$ws.Range("L132").Value = 9241.8948
$ws.Range("M132").Value = -5220.268999999999
$ws.Range("N132").Value = -14301.8948

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 797.75
$ws.Range("I93").Value = 609
$ws.Range("K93").Value = 609
$ws.Range("M93").Value = 639
$ws.Range("H132").Value = 5995.1
$ws.Range("I132").Value = 5455.8887
$ws.Range("J132").Value = 6436.273
$ws.Range("K132").Value = 16367.6661
$ws.Range("L132").Value = 19308.819
$ws.Range("M132").Value = -13837.6661
$ws.Range("N132").Value = -24368.819
$ws.Range("H136").Value = 5256.346
$ws.Range("I136").Value = 2761.5334
$ws.Range("J136").Value = 8658.362999999999
$ws.Range("K136").Value = 8284.600199999999
$ws.Range("L136").Value = 25975.089
$ws.Range("M136").Value = -5734.600199999999
$ws.Range("N136").Value = -31075.089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 25003612
$ws.Range("I132").Value = 45458020
$ws.Range("J132").Value = 3782.2222
$ws.Range("K132").Value = 136374060
$ws.Range("L132").Value = 11346.6666
$ws.Range("M132").Value = -136371530
$ws.Range("N132").Value = -16406.6666
$ws.Range("H136").Value = 10449784
$ws.Range("I136").Value = 22290058
$ws.Range("J136").Value = 2482.7058
$ws.Range("K136").Value = 66870174
$ws.Range("L136").Value = 7448.117400000001
$ws.Range("M136").Value = -66867624
$ws.Range("N136").Value = -12548.1174
